# Power consumption workbook update:
# add MiniAMP + speakers measurements (rows 10-14), fix row 9 styles,
# and widen column A / merge column D+E width definitions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dash = [char]0x2013

# --- Row 9 existing row: C9 / D9 style tweak (drop the now-unused duplicate xf) ---
$ws.Range("C9").Value = 200
$ws.Range("D9").Value = 3.7

# --- New data rows 10-14 ---
$labels = @(
    "RPi Zero W, MiniAMP, UPS Hat $dash No Speaker, 100% Volume",
    "RPi Zero W, MiniAMP, UPS Hat $dash 1 Speaker, 50% Volume, Over the Rainbow",
    "RPi Zero W, MiniAMP, UPS Hat $dash 1 Speaker, 50% Volume, 120Hz Sine",
    "RPi Zero W, MiniAMP, UPS Hat $dash 2 Speakers, 100% Volume, O. t. Rainbow, AVG!!",
    "RPi Zero W, MiniAMP, UPS Hat $dash 2 Speakers, 100% Volume, 120Hz Sine"
)
$playingMa = @(240, 250, 450, 700, 1260)

for ($i = 0; $i -lt 5; $i++) {
    $r = 10 + $i
    $ws.Range("A$r").Value = $labels[$i]
    $ws.Range("B$r").Value = 180
    $ws.Range("C$r").Value = $playingMa[$i]
    $ws.Range("D$r").Value = 3.7
    $ws.Range("E$r").Formula = "=B$r/1000*`$D$r"
    $ws.Range("F$r").Formula = "=C$r/1000*`$D$r"
    $ws.Range("G$r").Formula = "=F$r*(1+G`$4/100)"
    $ws.Range("H$r").Formula = "=((H`$4/1000*3.7)/`$G$r)"
    $ws.Range("I$r").Formula = "=((I`$4/1000*3.7)/`$G$r)"
    $ws.Range("J$r").Formula = "=((J`$4/1000*3.7)/`$G$r)"
    $ws.Range("K$r").Formula = "=((K`$4/1000*3.7)/`$G$r)"
}

# --- Re-select the cell the author left active ---
$ws.Range("A14").Select()

# --- Column A got wider to fit the longer device/setup descriptions ---
$ws.Columns.Item(1).ColumnWidth = 69.676666666666666

Write-Output "done"
